$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header order after removing "Kategorie" and "ID hráče" columns:
$ws.Range("A2").Value = "Jméno"
$ws.Range("B2").Value = "Přijmení"
$ws.Range("C2").Value = "Gender"
$ws.Range("D2").Value = "Klub"
$ws.Range("E2").Value = "Nasazení dle žebříčku"

# Clear the now-unused trailing columns but keep a centered style without border
$ws.Range("F2").ClearContents()
$ws.Range("G2").ClearContents()

$ws.Range("F2:G2").Borders.LineStyle = -4142
$ws.Range("F2:G2").HorizontalAlignment = -4108

$ws.Range("E2").Select()
